$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Applied Language"
$ws.Range("B1").Value = "Expected Data"

$ws.Range("B1").Select()
